$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns (zh-cn, de-de) and Latest Handoff Date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-26-17 12:26:42"

# zh-cn sheet: Status column and Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-17 12:26:39"

# de-de sheet: Status column and Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-17 12:26:42"
